$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 400
$ws.Range("I2").Value = 133
$ws.Range("J2").Value = 800.5
$ws.Range("K2").Value = 133
$ws.Range("L2").Value = 800.5
$ws.Range("M2").Value = -20
$ws.Range("N2").Value = -1026.5
$ws.Range("H41").Value = 1512.4
$ws.Range("I41").Value = 857.4286
$ws.Range("J41").Value = 1865.0769
$ws.Range("K41").Value = 857.4286
$ws.Range("L41").Value = 1865.0769
$ws.Range("M41").Value = -417.4286
$ws.Range("N41").Value = -2745.0769
$ws.Range("H43").Value = 1288299.9
$ws.Range("I43").Value = 1545159.9
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 1545159.9
$ws.Range("L43").Value = 4000
$ws.Range("M43").Value = -1545090.9
$ws.Range("N43").Value = -4138
$ws.Range("H76").Value = 3617.818
$ws.Range("I76").Value = 3527.7144
$ws.Range("J76").Value = 3775.5
$ws.Range("K76").Value = 3527.7144
$ws.Range("L76").Value = 3775.5
$ws.Range("M76").Value = -3212.7144
$ws.Range("N76").Value = -4405.5
$ws.Range("H79").Value = 3617.818
$ws.Range("I79").Value = 3527.7144
$ws.Range("J79").Value = 3775.5
$ws.Range("K79").Value = 3527.7144
$ws.Range("L79").Value = 3775.5
$ws.Range("M79").Value = -2435.7144
$ws.Range("N79").Value = -5959.5
$ws.Range("H86").Value = 4098283.5
$ws.Range("J86").Value = 4551.25
$ws.Range("L86").Value = 4551.25
$ws.Range("N86").Value = -6797.25
$ws.Range("H89").Value = 4098283.5
$ws.Range("J89").Value = 4551.25
$ws.Range("L89").Value = 22756.25
$ws.Range("N89").Value = -33988.25
$ws.Range("H99").Value = 2160.6
$ws.Range("I99").Value = 2475.75
$ws.Range("K99").Value = 7427.25
$ws.Range("M99").Value = -5929.25
$ws.Range("H107").Value = 312
$ws.Range("I107").Value = 356.375
$ws.Range("J107").Value = 261.2857
$ws.Range("K107").Value = 356.375
$ws.Range("L107").Value = 261.2857
$ws.Range("M107").Value = 1563.625
$ws.Range("N107").Value = -4101.2857
$ws.Range("H115").Value = 781.6667
$ws.Range("I115").Value = 781.6667
$ws.Range("K115").Value = 2345.0001
$ws.Range("M115").Value = -778.0001000000002
$ws.Range("H118").Value = 879.5
$ws.Range("I118").Value = 939.3333
$ws.Range("K118").Value = 2817.9999
$ws.Range("M118").Value = -1160.9999
$ws.Range("H127").Value = 2651.537
$ws.Range("I127").Value = 1038.5625
$ws.Range("K127").Value = 3115.6875
$ws.Range("M127").Value = 1844.3125
$ws.Range("H129").Value = 1825.4482
$ws.Range("I129").Value = 558.2857
$ws.Range("K129").Value = 1674.8571
$ws.Range("M129").Value = 3325.1429
$ws.Range("H137").Value = 5876
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 5876
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 17628
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -22728
$ws.Range("H138").Value = 5705.0894
$ws.Range("I138").Value = 1527.75
$ws.Range("J138").Value = 7376.025
$ws.Range("K138").Value = 4583.25
$ws.Range("L138").Value = 22128.075
$ws.Range("M138").Value = 556.75
$ws.Range("N138").Value = -32408.075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3922.8604
$ws.Range("I32").Value = 1810.8334
$ws.Range("K32").Value = 1810.8334
$ws.Range("M32").Value = -1523.8334
$ws.Range("H63").Value = 3162.5557
$ws.Range("J63").Value = 3227.75
$ws.Range("L63").Value = 3227.75
$ws.Range("N63").Value = -4599.75
$ws.Range("H66").Value = 3162.5557
$ws.Range("J66").Value = 3227.75
$ws.Range("L66").Value = 16138.75
$ws.Range("N66").Value = -23002.75
$ws.Range("H132").Value = 24359.172
$ws.Range("I132").Value = 28181.809
$ws.Range("K132").Value = 84545.427
$ws.Range("M132").Value = -82015.427

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 415724.03
$ws.Range("I94").Value = 856426.75
$ws.Range("K94").Value = 856426.75
$ws.Range("M94").Value = -855975.75
$ws.Range("H107").Value = 1124.0333
$ws.Range("I107").Value = 1093.8667
$ws.Range("K107").Value = 1093.8667
$ws.Range("M107").Value = 826.1333
$ws.Range("H134").Value = 2645.4736
$ws.Range("I134").Value = 2149.4517
$ws.Range("K134").Value = 6448.355100000001
$ws.Range("M134").Value = -3913.355100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2327.2856
$ws.Range("I31").Value = 1048.5
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 1048.5
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -753.5
$ws.Range("N31").Value = -10590
$ws.Range("H34").Value = 2327.2856
$ws.Range("I34").Value = 1048.5
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 1048.5
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -846.5
$ws.Range("N34").Value = -10404
$ws.Range("H99").Value = 11794.25
$ws.Range("I99").Value = 20285
$ws.Range("K99").Value = 20285
$ws.Range("M99").Value = -18787
$ws.Range("H126").Value = 11794.25
$ws.Range("I126").Value = 20285
$ws.Range("K126").Value = 60855
$ws.Range("M126").Value = -58385
$ws.Range("H134").Value = 2098.682
$ws.Range("I134").Value = 2051
$ws.Range("K134").Value = 6153
$ws.Range("M134").Value = -3618
$ws.Range("H140").Value = 100780
$ws.Range("J140").Value = 100780
$ws.Range("L140").Value = 100780
$ws.Range("N140").Value = -111140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 20998.6
$ws.Range("J123").Value = 24997.666
$ws.Range("L123").Value = 74992.99800000001
$ws.Range("N123").Value = -79892.99800000001
$ws.Range("H131").Value = 7248830
$ws.Range("I131").Value = 23810514
$ws.Range("J131").Value = 5378962.5
$ws.Range("K131").Value = 71431542
$ws.Range("L131").Value = 16136887.5
$ws.Range("M131").Value = -71426502
$ws.Range("N131").Value = -16146967.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 892
$ws.Range("I97").Value = 892
$ws.Range("K97").Value = 892
$ws.Range("M97").Value = -396
$ws.Range("H102").Value = 6908.154
$ws.Range("I102").Value = 6529.619
$ws.Range("K102").Value = 6529.619
$ws.Range("M102").Value = -4907.619
$ws.Range("H122").Value = 1001962.06
$ws.Range("I122").Value = 1835169.1
$ws.Range("K122").Value = 5505507.300000001
$ws.Range("M122").Value = -5503057.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8335093
$ws.Range("I16").Value = 11112374
$ws.Range("K16").Value = 11112374
$ws.Range("M16").Value = -11112204
$ws.Range("H61").Value = 3892.9
$ws.Range("I61").Value = 4532.857
$ws.Range("K61").Value = 4532.857
$ws.Range("M61").Value = -4330.857
$ws.Range("H93").Value = 903.94446
$ws.Range("I93").Value = 804.73334
$ws.Range("K93").Value = 804.73334
$ws.Range("M93").Value = 443.26666
$ws.Range("H113").Value = 3892.9
$ws.Range("I113").Value = 4532.857
$ws.Range("K113").Value = 4532.857
$ws.Range("M113").Value = -2362.857
$ws.Range("H132").Value = 3722.8533
$ws.Range("I132").Value = 2817.1035
$ws.Range("K132").Value = 8451.3105
$ws.Range("M132").Value = -5921.3105

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2203104
$ws.Range("I81").Value = 2086902.8
$ws.Range("J81").Value = 2332216.2
$ws.Range("K81").Value = 4173805.6
$ws.Range("L81").Value = 4664432.4
$ws.Range("M81").Value = -4172744.6
$ws.Range("N81").Value = -4666554.4
$ws.Range("H84").Value = 2203104
$ws.Range("I84").Value = 2086902.8
$ws.Range("J84").Value = 2332216.2
$ws.Range("K84").Value = 20869028
$ws.Range("L84").Value = 23322162
$ws.Range("M84").Value = -20863724
$ws.Range("N84").Value = -23332770
$ws.Range("H96").Value = 5526.625
$ws.Range("I96").Value = 4003
$ws.Range("K96").Value = 4003
$ws.Range("M96").Value = -2630
$ws.Range("H132").Value = 18523236
$ws.Range("I132").Value = 1933.0952
$ws.Range("J132").Value = 83347800
$ws.Range("K132").Value = 5799.2856
$ws.Range("L132").Value = 250043400
$ws.Range("M132").Value = -3269.2856
$ws.Range("N132").Value = -250048460
